$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 79, shifting the existing data (rows 79-105)
# down to rows 80-106.
$ws.Rows.Item(79).Insert()

# Populate the newly inserted row 79 with the new weekly price record.
$ws.Cells.Item(79, 1).Value = 4
$ws.Cells.Item(79, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(79, 3).Value = "Los Lagos"
$ws.Cells.Item(79, 4).Value = 44463
$ws.Cells.Item(79, 5).Value = 10
$ws.Cells.Item(79, 6).Value = 100112009
$ws.Cells.Item(79, 7).Value = "Acelga"
$ws.Cells.Item(79, 8).Value = "Sin especificar"
$ws.Cells.Item(79, 9).Value = "Primera"
$ws.Cells.Item(79, 10).Value = 200
$ws.Cells.Item(79, 11).Value = 4000
$ws.Cells.Item(79, 12).Value = 4000
$ws.Cells.Item(79, 13).Value = 4000
$ws.Cells.Item(79, 14).Value = "`$/docena de atados (4 kilos)"
$ws.Cells.Item(79, 15).Value = "Región del Maule"
$ws.Cells.Item(79, 16).Value = 1000
$ws.Cells.Item(79, 17).Value = 4
$ws.Cells.Item(79, 18).Value = "Hortaliza"
